# Updates as of 8th April 2020
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data row for 8 April 2020 (row 26)
$row = 26

$ws.Cells.Item($row, 1).Value = 43929                          # A26 - Date (8-Apr-2020 serial)
$ws.Cells.Item($row, 1).NumberFormat = "d-mmm-yy"               # match the format used by the rest of column A
$ws.Cells.Item($row, 1).HorizontalAlignment = -4108              # xlCenter
$ws.Cells.Item($row, 2).Value = 7                              # B26 - New Cases
$ws.Cells.Item($row, 3).Value = 305                             # C26 - Tested

# Set the new text values in the same order the new shared strings were
# originally authored in (Case Type, then Travelled From, then County) so
# the resulting shared-string table indices line up with the source edit.
$ws.Cells.Item($row, 7).Value = "Community(4), Imported(3)"      # G26 - Case Type
$ws.Cells.Item($row, 4).Value = "Congo, United States, UK(2)"   # D26 - Travelled From
$ws.Cells.Item($row, 5).Value = "Nairobi, Mombasa, Uasin Ngishu" # E26 - County

$ws.Cells.Item($row, 6).Value = 179                              # F26 - Aggregation
$ws.Cells.Item($row, 8).Value = 0                                 # H26 - Recover
$ws.Cells.Item($row, 9).Value = 0                                 # I26 - Death
$ws.Cells.Item($row, 11).Value = "Mercy"                          # K26 - Info Giver

# Update the sheet view to reflect the newly selected/visible cells
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L26").Select() | Out-Null
